$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "81 x 60" + [char]11 + "  6    0" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(1,2).Range.Text = "20 x 61" + [char]11 + "  6    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "0|    |"
$t.Cell(1,3).Range.Text = "79 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"
$t.Cell(2,1).Range.Text = "99 x 32" + [char]11 + "  3    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "9|    |"
$t.Cell(2,2).Range.Text = "40 x 53" + [char]11 + "  5    3" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"
$t.Cell(2,3).Range.Text = "43 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
$t.Cell(3,1).Range.Text = "42 x 66" + [char]11 + "  6    6" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "2|    |"
$t.Cell(3,2).Range.Text = "18 x 12" + [char]11 + "  1    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "8|    |"
$t.Cell(3,3).Range.Text = "58 x 15" + [char]11 + "  1    5" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "8|    |"
$t.Cell(4,1).Range.Text = "66 x 61" + [char]11 + "  6    1" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(4,2).Range.Text = "35 x 40" + [char]11 + "  4    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
$t.Cell(4,3).Range.Text = "66 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(5,1).Range.Text = "51 x 10" + [char]11 + "  1    0" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "1|    |"
$t.Cell(5,2).Range.Text = "71 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "1|    |"
$t.Cell(5,3).Range.Text = "33 x 20" + [char]11 + "  2    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
